$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Range("B1").Value = "first feed number of pellets left"
$ws.Range("C1").Value = "second feed number of pellets left"
$ws.Range("D1").Value = "total feed pellets fed"

# Copy header style (bold) from B1 to C1 and D1
$ws.Range("B1").Copy() | Out-Null
$ws.Range("C1:D1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Update data rows
$dates = @("26 Nov", "27 Nov", "28 Nov", "29 Nov", "30 Nov", "01 Dec", "02 Dec")
$firstFeed = @(0, 33, 29, 30, 32, 34, 31)
$secondFeed = @(0, 42, 41, 38, 35, 35, 35)
$total = @(0, 75, 70, 68, 67, 69, 66)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $dates[$i]
    $ws.Cells.Item($row, 2).Value = $firstFeed[$i]
    $ws.Cells.Item($row, 3).Value = $secondFeed[$i]
    $ws.Cells.Item($row, 4).Value = $total[$i]
}
